$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-06 Tuesday" "2024-02-07 Wednesday"

Replace-Text "48×38=" "43×45="
Replace-Text "75×36=" "82×69="
Replace-Text "11×18=" "41×60="
Replace-Text "75×44=" "15×36="
Replace-Text "27×98=" "23×36="

Replace-Text "73×30=" "86×80="
Replace-Text "80×69=" "70×58="
Replace-Text "94×42=" "32×57="
Replace-Text "82×34=" "56×40="
Replace-Text "76×99=" "94×25="

Replace-Text "67×64=" "79×66="
Replace-Text "72×27=" "69×71="
Replace-Text "39×43=" "58×40="
Replace-Text "59×45=" "24×68="
Replace-Text "23×65=" "11×23="

Replace-Text "47×67=" "99×91="
Replace-Text "55×84=" "37×31="
Replace-Text "40×88=" "22×25="
Replace-Text "93×77=" "41×43="
Replace-Text "19×17=" "37×39="

Replace-Text "91×38=" "41×38="
Replace-Text "61×97=" "61×18="
Replace-Text "55×93=" "18×50="
Replace-Text "28×95=" "26×32="
Replace-Text "86×70=" "64×60="
